$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match latest scrape.
# NumberFormat "@" forces text entry (preventing Excel from auto-converting
# numeric-looking strings into floating point numbers), then ClearFormats()
# removes the temporary format so the cell keeps its original (default) style.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "56.644.01"
Set-TextValue "E2" "  +10.62%  "
Set-TextValue "D3" "3.255.28"
Set-TextValue "E3" "  +6.24%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "398.11"
Set-TextValue "E5" "  +2.62%  "
Set-TextValue "D6" "111.13"
Set-TextValue "E6" "  +8.88%  "
Set-TextValue "E7" "  +4.58%  "
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "D9" "0.623"
Set-TextValue "E9" "  +7.60%  "
Set-TextValue "D10" "39.46"
Set-TextValue "E10" "  +7.33%  "
Set-TextValue "D11" "0.0955"
Set-TextValue "E11" "  +12.56%  "
Set-TextValue "E12" "  +2.26%  "
Set-TextValue "D13" "3.749.30"
Set-TextValue "E13" "  +5.74%  "
Set-TextValue "E14" "  +5.73%  "
Set-TextValue "D15" "8.14"
Set-TextValue "E15" "  +5.93%  "
Set-TextValue "D16" "3.250.41"
Set-TextValue "E16" "  +6.24%  "
Set-TextValue "E17" "  +5.79%  "
Set-TextValue "D18" "11.12"
Set-TextValue "E18" "  +4.12%  "
Set-TextValue "D19" "56.586.69"
Set-TextValue "E19" "  +10.49%  "
Set-TextValue "E20" "  +4.21%  "
Set-TextValue "D21" "0.0000104"
Set-TextValue "E21" "  +9.43%  "
Set-TextValue "D22" "13.04"
Set-TextValue "E22" "  +6.48%  "
Set-TextValue "D23" "300.08"
Set-TextValue "E23" "  +13.53%  "
Set-TextValue "D24" "75.53"
Set-TextValue "E24" "  +8.49%  "
Set-TextValue "E25" "  +3.89%  "
Set-TextValue "D26" "8.16"
Set-TextValue "E26" "  +3.17%  "
Set-TextValue "D27" "28.35"
Set-TextValue "E27" "  +4.77%  "
Set-TextValue "E28" "  +4.17%  "
Set-TextValue "E29" "  +1.48%  "
Set-TextValue "D30" "0.170"
Set-TextValue "E30" "  +4.78%  "
Set-TextValue "E31" "  -0.06%  "
Set-TextValue "E32" "  +6.80%  "
Set-TextValue "D33" "11.10"
Set-TextValue "E33" "  +6.59%  "
Set-TextValue "D34" "36.79"
Set-TextValue "E34" "  +3.18%  "
Set-TextValue "E35" "  +2.82%  "
Set-TextValue "E36" "  +2.12%  "
Set-TextValue "D37" "51.69"
Set-TextValue "E37" "  +3.25%  "
Set-TextValue "E38" "  +26.87%  "
Set-TextValue "E39" "  +5.73%  "
Set-TextValue "D40" "1.00"
Set-TextValue "E40" "  +0.07%  "
Set-TextValue "D41" "17.62"
Set-TextValue "D42" "134.37"
Set-TextValue "E42" "  +3.02%  "
Set-TextValue "E43" "  +5.94%  "
Set-TextValue "E44" "  +4.78%  "
Set-TextValue "D45" "3.99"
Set-TextValue "E45" "  +6.64%  "
Set-TextValue "E46" "  -3.51%  "
Set-TextValue "D47" "22.29"
Set-TextValue "E47" "  +2.95%  "
Set-TextValue "D48" "2.21"
Set-TextValue "E48" "  +57.11%  "
Set-TextValue "D49" "2.149.95"
Set-TextValue "E49" "  +4.30%  "
Set-TextValue "E50" "  +1.11%  "
Set-TextValue "D51" "2.41"
Set-TextValue "E51" "  -3.87%  "
